# Refresh the cryptos table: update Price (D) and Volume(1h) (E) columns
# for each coin row whose market data moved since the last snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.536.62"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "2.633.87"
$ws.Range("E3").Value = "  -1.63%  "

$ws.Range("D5").Value = "'595.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D6").Value = "'168.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -2.43%  "

$ws.Range("D9").Value = "2.633.66"
$ws.Range("E9").Value = "  -1.62%  "

$ws.Range("E10").Value = "  -3.14%  "

$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("E13").Value = "  -0.15%  "

$ws.Range("D14").Value = "'27.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.81%  "

$ws.Range("D15").Value = "3.113.21"
$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("D17").Value = "67.423.83"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "2.612.99"
$ws.Range("E18").Value = "  -2.41%  "

$ws.Range("E19").Value = "  +1.72%  "

$ws.Range("E20").Value = "  +3.78%  "

$ws.Range("D21").Value = "'357.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.90%  "

$ws.Range("D22").Value = "'4.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.91%  "

$ws.Range("E24").Value = "  -4.65%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "'10.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.68%  "

$ws.Range("D27").Value = "'69.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("D31").Value = "'548.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.72%  "

$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("E34").Value = "  -2.09%  "

$ws.Range("E35").Value = "  +4.51%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  -4.03%  "

$ws.Range("D38").Value = "'157.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("D39").Value = "'19.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.66%  "

$ws.Range("D40").Value = "'0.365"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.33%  "

$ws.Range("E41").Value = "  -0.62%  "

$ws.Range("D42").Value = "'18.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.99%  "

$ws.Range("D43").Value = "'5.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.79%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("E45").Value = "  -3.95%  "

$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").Value = "'153.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").Value = "'0.580"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("D49").Value = "'3.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.50%  "

$ws.Range("D50").Value = "'1.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "

$ws.Range("E51").Value = "  -1.09%  "

Write-Output "Updated cryptos list"
